# Fix: when all observations in an aggregating group are NA, sum(na.rm=T)
# previously returned 0; now it should return NA. Update the "Value" column
# (column G) on the "Data" sheet from "0" to "NaN" for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$rows = @(3,6,10,11,12,13,17,19,20,23,27,28,31,286,292,294,298,300,301,302,310,311,569,575,577,581,583,584,585,593,594,860,866,867,868,876,877)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Value"
    if ($cell.Value2 -eq "0") {
        $cell.Value2 = "NaN"
    }
}
